# Work on import pointsourceactivities:
# append two new EMEP emission-factor rows (NFR 1.A.4.b.i / 1.A.4.b.ii) to
# the short EMEP emission factors table, and touch a couple of cosmetic
# workbook/view settings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cosmetic workbook/view tweaks -----------------------------------
$wb.Windows.Item(1).TabRatio = 0.5
$excel.ActiveWindow.Zoom = 100
$ws.StandardWidth = 11.55078125

# --- new row 12: 1.A.4.b.i / Residential plants / PCB -----------------
$ws.Cells.Item(12, 1).Value = "1.A.4.b.i"
$ws.Cells.Item(12, 2).Value = "Residential plants"
$ws.Cells.Item(12, 3).Value = "Table_3-19"
$ws.Cells.Item(12, 4).Value = "Tier 2 emission factor"
$ws.Cells.Item(12, 5).Value = "Advanced coal combustion techniques <1MWth - Advanced stove"
$ws.Cells.Item(12, 6).Value = "Coal Fuels"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "PCB"
$ws.Cells.Item(12, 10).Value = 170
$ws.Cells.Item(12, 11).Value = "µg/GJ"
$ws.Cells.Item(12, 12).Value = 85
$ws.Cells.Item(12, 13).Value = 260
$ws.Cells.Item(12, 14).Value = "Kakareka et al. (2004)"

# --- new row 13: 1.A.4.b.ii / Household and gardening (mobile) / NMVOC -
$ws.Cells.Item(13, 1).Value = "1.A.4.b.ii"
$ws.Cells.Item(13, 2).Value = "Household and gardening (mobile)"
$ws.Cells.Item(13, 3).Value = "Table_3-1_05"
$ws.Cells.Item(13, 4).Value = "Tier 1 emission factor"
$ws.Cells.Item(13, 5).Value = "NA"
$ws.Cells.Item(13, 6).Value = "LPG"
$ws.Cells.Item(13, 8).Value = "NA"
$ws.Cells.Item(13, 9).Value = "NMVOC"
$ws.Cells.Item(13, 10).Value = 6720
$ws.Cells.Item(13, 11).Value = "g/tonnes fuel"
$ws.Cells.Item(13, 14).Value = "Winther 2016"

# match the rest of the sheet's styling (Arial 10) for the new rows
$ws.Range("A12:W13").Font.Name = "Arial"
$ws.Range("A12:W13").Font.Size = 10

# select the newly added rows, like a user reviewing the freshly pasted data
$ws.Rows("12:13").Select() | Out-Null
